# Add "FINAL GROUP PROJECT - " in front of the existing "TEXT DOCUMENT"
# title run, and move the auto-managed "_GoBack" bookmark (Word drops
# this at the location of the most recent edit) from its old location
# -- just before the "Include a brief statement..." run -- to sit right
# after the freshly typed text, matching where Word would leave it after
# this edit.

$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark from its old position ---------
# Word always has at most one "_GoBack" bookmark; it is hidden from the
# regular Bookmarks collection/count but can still be addressed by name.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Type the new heading text in front of "TEXT DOCUMENT" -------
$titlePara = $d.Paragraphs(1).Range
$insertionPoint = $d.Range($titlePara.Start, $titlePara.Start)

$newText = "FINAL GROUP PROJECT - "
$insertionPoint.InsertBefore($newText)

# --- 3. Drop "_GoBack" back in, now marking this latest edit --------
$afterNewText = $titlePara.Start + $newText.Length
$newMark = $d.Range($afterNewText, $afterNewText)
$d.Bookmarks.Add("_GoBack", $newMark)
